$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.464.16"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.640.55"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9978"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.78"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3770"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.70"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08191"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.36"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.528"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.340"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "1.639.84"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.28"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06971"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.739"
$ws.Range("E20").Value = "  +4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.46"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Value = "23.467.30"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.519"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.126"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.55"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.182"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.20"
$ws.Range("D31").Value = "1.825.77"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.762"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.091"
$ws.Range("E33").Value = "  +7.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.48"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.032"
$ws.Range("E35").Value = "  -10.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02765"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2495"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08775"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.019"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06991"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7006"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.323"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.62"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6457"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.327"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9983"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.957"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07962"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.54"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.180"
$ws.Range("E51").Value = "  -1.92%  "
